$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: last-updated timestamp banner
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 1 de Abril de 2020 a las 13:20"

# Row 3: column headers (unchanged, kept for completeness)
$ws.Cells.Item(3, 1).Value = "Ciudad"
$ws.Cells.Item(3, 2).Value = "Casos totales"
$ws.Cells.Item(3, 3).Value = "Casos activos"
$ws.Cells.Item(3, 4).Value = "Recuperados"
$ws.Cells.Item(3, 5).Value = "Muertes"

# Rows 4-64: per-province data (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
$ws.Cells.Item(4, 1).Value = "Madrid"
$ws.Cells.Item(4, 2).Value = 29840
$ws.Cells.Item(4, 3).Value = 10827
$ws.Cells.Item(4, 4).Value = 15148
$ws.Cells.Item(4, 5).Value = 3865
$ws.Cells.Item(5, 1).Value = "Cataluña"
$ws.Cells.Item(5, 2).Value = 19991
$ws.Cells.Item(5, 3).Value = 5701
$ws.Cells.Item(5, 4).Value = 12441
$ws.Cells.Item(5, 5).Value = 1849
$ws.Cells.Item(6, 1).Value = "Bizkaia/Vizcaya"
$ws.Cells.Item(6, 2).Value = 3382
$ws.Cells.Item(6, 3).Value = 2165
$ws.Cells.Item(6, 4).Value = 2289
$ws.Cells.Item(6, 5).Value = 168
$ws.Cells.Item(7, 1).Value = "Valencia/Valencia"
$ws.Cells.Item(7, 2).Value = 3089
$ws.Cells.Item(7, 3).Value = 169
$ws.Cells.Item(7, 4).Value = 2734
$ws.Cells.Item(7, 5).Value = 186
$ws.Cells.Item(8, 1).Value = "Castilla-La Mancha"
$ws.Cells.Item(8, 2).Value = 2780
$ws.Cells.Item(8, 3).Value = 71
$ws.Cells.Item(8, 4).Value = 2446
$ws.Cells.Item(8, 5).Value = 263
$ws.Cells.Item(9, 1).Value = "Ciudad Real"
$ws.Cells.Item(9, 2).Value = 2471
$ws.Cells.Item(9, 3).Value = 397
$ws.Cells.Item(9, 4).Value = 2076
$ws.Cells.Item(9, 5).Value = 245
$ws.Cells.Item(10, 1).Value = "Navarra"
$ws.Cells.Item(10, 2).Value = 2305
$ws.Cells.Item(10, 3).Value = 192
$ws.Cells.Item(10, 4).Value = 2000
$ws.Cells.Item(10, 5).Value = 113
$ws.Cells.Item(11, 1).Value = "Araba/Alava"
$ws.Cells.Item(11, 2).Value = 2250
$ws.Cells.Item(11, 3).Value = 2165
$ws.Cells.Item(11, 4).Value = 1376
$ws.Cells.Item(11, 5).Value = 149
$ws.Cells.Item(12, 1).Value = "Alacant/Alicante"
$ws.Cells.Item(12, 2).Value = 2173
$ws.Cells.Item(12, 3).Value = 62
$ws.Cells.Item(12, 4).Value = 1944
$ws.Cells.Item(12, 5).Value = 167
$ws.Cells.Item(13, 1).Value = "La Rioja"
$ws.Cells.Item(13, 2).Value = 1960
$ws.Cells.Item(13, 3).Value = 569
$ws.Cells.Item(13, 4).Value = 1301
$ws.Cells.Item(13, 5).Value = 90
$ws.Cells.Item(14, 1).Value = "Albacete"
$ws.Cells.Item(14, 2).Value = 1933
$ws.Cells.Item(14, 3).Value = 397
$ws.Cells.Item(14, 4).Value = 1678
$ws.Cells.Item(14, 5).Value = 156
$ws.Cells.Item(15, 1).Value = "Zaragoza"
$ws.Cells.Item(15, 2).Value = 1902
$ws.Cells.Item(15, 3).Value = 208
$ws.Cells.Item(15, 4).Value = 1564
$ws.Cells.Item(15, 5).Value = 130
$ws.Cells.Item(16, 1).Value = "A Coruña"
$ws.Cells.Item(16, 2).Value = 1854
$ws.Cells.Item(16, 3).Value = 259
$ws.Cells.Item(16, 4).Value = 1706
$ws.Cells.Item(16, 5).Value = 58
$ws.Cells.Item(17, 1).Value = "Toledo"
$ws.Cells.Item(17, 2).Value = 1593
$ws.Cells.Item(17, 3).Value = 397
$ws.Cells.Item(17, 4).Value = 1298
$ws.Cells.Item(17, 5).Value = 205
$ws.Cells.Item(18, 1).Value = "Pontevedra"
$ws.Cells.Item(18, 2).Value = 1452
$ws.Cells.Item(18, 3).Value = 259
$ws.Cells.Item(18, 4).Value = 1337
$ws.Cells.Item(18, 5).Value = 30
$ws.Cells.Item(19, 1).Value = "Malaga"
$ws.Cells.Item(19, 2).Value = 1349
$ws.Cells.Item(19, 3).Value = 83
$ws.Cells.Item(19, 4).Value = 1196
$ws.Cells.Item(19, 5).Value = 70
$ws.Cells.Item(20, 1).Value = "Asturias"
$ws.Cells.Item(20, 2).Value = 1322
$ws.Cells.Item(20, 3).Value = 109
$ws.Cells.Item(20, 4).Value = 1150
$ws.Cells.Item(20, 5).Value = 63
$ws.Cells.Item(21, 1).Value = "Salamanca"
$ws.Cells.Item(21, 2).Value = 1316
$ws.Cells.Item(21, 3).Value = 235
$ws.Cells.Item(21, 4).Value = 946
$ws.Cells.Item(21, 5).Value = 135
$ws.Cells.Item(22, 1).Value = "Tenerife"
$ws.Cells.Item(22, 2).Value = 1262
$ws.Cells.Item(22, 3).Value = 30
$ws.Cells.Item(22, 4).Value = 1056
$ws.Cells.Item(22, 5).Value = 36
$ws.Cells.Item(23, 1).Value = "Gran Canaria"
$ws.Cells.Item(23, 2).Value = 1262
$ws.Cells.Item(23, 3).Value = 57
$ws.Cells.Item(23, 4).Value = 342
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(24, 1).Value = "La Palma"
$ws.Cells.Item(24, 2).Value = 1262
$ws.Cells.Item(24, 3).Value = 57
$ws.Cells.Item(24, 4).Value = 57
$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(25, 1).Value = "Lanzarote"
$ws.Cells.Item(25, 2).Value = 1262
$ws.Cells.Item(25, 3).Value = 57
$ws.Cells.Item(25, 4).Value = 45
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(26, 1).Value = "Fuerteventura"
$ws.Cells.Item(26, 2).Value = 1262
$ws.Cells.Item(26, 3).Value = 57
$ws.Cells.Item(26, 4).Value = 32
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(27, 1).Value = "La Gomera"
$ws.Cells.Item(27, 2).Value = 1262
$ws.Cells.Item(27, 3).Value = 57
$ws.Cells.Item(27, 4).Value = 7
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(28, 1).Value = "El Hierro"
$ws.Cells.Item(28, 2).Value = 1262
$ws.Cells.Item(28, 3).Value = 57
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 1).Value = "Cantabria"
$ws.Cells.Item(29, 2).Value = 1213
$ws.Cells.Item(29, 3).Value = 43
$ws.Cells.Item(29, 4).Value = 1116
$ws.Cells.Item(29, 5).Value = 54
$ws.Cells.Item(30, 1).Value = "Gipuzkoa/Guipuzcoa"
$ws.Cells.Item(30, 2).Value = 1206
$ws.Cells.Item(30, 3).Value = 2165
$ws.Cells.Item(30, 4).Value = 639
$ws.Cells.Item(30, 5).Value = 52
$ws.Cells.Item(31, 1).Value = "Sevilla"
$ws.Cells.Item(31, 2).Value = 1119
$ws.Cells.Item(31, 3).Value = 17
$ws.Cells.Item(31, 4).Value = 1066
$ws.Cells.Item(31, 5).Value = 36
$ws.Cells.Item(32, 1).Value = "Valladolid"
$ws.Cells.Item(32, 2).Value = 1109
$ws.Cells.Item(32, 3).Value = 262
$ws.Cells.Item(32, 4).Value = 758
$ws.Cells.Item(32, 5).Value = 89
$ws.Cells.Item(33, 1).Value = "Caceres"
$ws.Cells.Item(33, 2).Value = 1093
$ws.Cells.Item(33, 3).Value = 31
$ws.Cells.Item(33, 4).Value = 932
$ws.Cells.Item(33, 5).Value = 130
$ws.Cells.Item(34, 1).Value = "Granada"
$ws.Cells.Item(34, 2).Value = 1061
$ws.Cells.Item(34, 3).Value = 15
$ws.Cells.Item(34, 4).Value = 979
$ws.Cells.Item(34, 5).Value = 67
$ws.Cells.Item(35, 1).Value = "Murcia"
$ws.Cells.Item(35, 2).Value = 1041
$ws.Cells.Item(35, 3).Value = 43
$ws.Cells.Item(35, 4).Value = 961
$ws.Cells.Item(35, 5).Value = 37
$ws.Cells.Item(36, 1).Value = "Leon"
$ws.Cells.Item(36, 2).Value = 1007
$ws.Cells.Item(36, 3).Value = 216
$ws.Cells.Item(36, 4).Value = 683
$ws.Cells.Item(36, 5).Value = 108
$ws.Cells.Item(37, 1).Value = "Aragon"
$ws.Cells.Item(37, 2).Value = 907
$ws.Cells.Item(37, 3).Value = 29
$ws.Cells.Item(37, 4).Value = 838
$ws.Cells.Item(37, 5).Value = 40
$ws.Cells.Item(38, 1).Value = "Burgos"
$ws.Cells.Item(38, 2).Value = 848
$ws.Cells.Item(38, 3).Value = 250
$ws.Cells.Item(38, 4).Value = 521
$ws.Cells.Item(38, 5).Value = 77
$ws.Cells.Item(39, 1).Value = "Segovia"
$ws.Cells.Item(39, 2).Value = 811
$ws.Cells.Item(39, 3).Value = 210
$ws.Cells.Item(39, 4).Value = 520
$ws.Cells.Item(39, 5).Value = 81
$ws.Cells.Item(40, 1).Value = "Guadalajara"
$ws.Cells.Item(40, 2).Value = 753
$ws.Cells.Item(40, 3).Value = 397
$ws.Cells.Item(40, 4).Value = 618
$ws.Cells.Item(40, 5).Value = 100
$ws.Cells.Item(41, 1).Value = "Cordoba"
$ws.Cells.Item(41, 2).Value = 661
$ws.Cells.Item(41, 3).Value = 4
$ws.Cells.Item(41, 4).Value = 642
$ws.Cells.Item(41, 5).Value = 15
$ws.Cells.Item(42, 1).Value = "Jaen"
$ws.Cells.Item(42, 2).Value = 661
$ws.Cells.Item(42, 3).Value = 17
$ws.Cells.Item(42, 4).Value = 618
$ws.Cells.Item(42, 5).Value = 26
$ws.Cells.Item(43, 1).Value = "Castello/Castellon"
$ws.Cells.Item(43, 2).Value = 660
$ws.Cells.Item(43, 3).Value = 9
$ws.Cells.Item(43, 4).Value = 609
$ws.Cells.Item(43, 5).Value = 42
$ws.Cells.Item(44, 1).Value = "Soria"
$ws.Cells.Item(44, 2).Value = 659
$ws.Cells.Item(44, 3).Value = 90
$ws.Cells.Item(44, 4).Value = 525
$ws.Cells.Item(44, 5).Value = 44
$ws.Cells.Item(45, 1).Value = "Ourense"
$ws.Cells.Item(45, 2).Value = 626
$ws.Cells.Item(45, 3).Value = 259
$ws.Cells.Item(45, 4).Value = 559
$ws.Cells.Item(45, 5).Value = 18
$ws.Cells.Item(46, 1).Value = "Badajoz"
$ws.Cells.Item(46, 2).Value = 586
$ws.Cells.Item(46, 3).Value = 82
$ws.Cells.Item(46, 4).Value = 482
$ws.Cells.Item(46, 5).Value = 22
$ws.Cells.Item(47, 1).Value = "Cadiz"
$ws.Cells.Item(47, 2).Value = 539
$ws.Cells.Item(47, 3).Value = 16
$ws.Cells.Item(47, 4).Value = 509
$ws.Cells.Item(47, 5).Value = 14
$ws.Cells.Item(48, 1).Value = "Avila"
$ws.Cells.Item(48, 2).Value = 512
$ws.Cells.Item(48, 3).Value = 132
$ws.Cells.Item(48, 4).Value = 321
$ws.Cells.Item(48, 5).Value = 59
$ws.Cells.Item(49, 1).Value = "Lugo"
$ws.Cells.Item(49, 2).Value = 500
$ws.Cells.Item(49, 3).Value = 259
$ws.Cells.Item(49, 4).Value = 456
$ws.Cells.Item(49, 5).Value = 9
$ws.Cells.Item(50, 1).Value = "Palencia"
$ws.Cells.Item(50, 2).Value = 359
$ws.Cells.Item(50, 3).Value = 52
$ws.Cells.Item(50, 4).Value = 284
$ws.Cells.Item(50, 5).Value = 23
$ws.Cells.Item(51, 1).Value = "Huesca"
$ws.Cells.Item(51, 2).Value = 317
$ws.Cells.Item(51, 3).Value = 28
$ws.Cells.Item(51, 4).Value = 273
$ws.Cells.Item(51, 5).Value = 16
$ws.Cells.Item(52, 1).Value = "Cuenca"
$ws.Cells.Item(52, 2).Value = 297
$ws.Cells.Item(52, 3).Value = 397
$ws.Cells.Item(52, 4).Value = 206
$ws.Cells.Item(52, 5).Value = 68
$ws.Cells.Item(53, 1).Value = "Teruel"
$ws.Cells.Item(53, 2).Value = 272
$ws.Cells.Item(53, 3).Value = 21
$ws.Cells.Item(53, 4).Value = 230
$ws.Cells.Item(53, 5).Value = 21
$ws.Cells.Item(54, 1).Value = "Almeria"
$ws.Cells.Item(54, 2).Value = 251
$ws.Cells.Item(54, 3).Value = 6
$ws.Cells.Item(54, 4).Value = 229
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(55, 1).Value = "Zamora"
$ws.Cells.Item(55, 2).Value = 226
$ws.Cells.Item(55, 3).Value = 51
$ws.Cells.Item(55, 4).Value = 150
$ws.Cells.Item(55, 5).Value = 25
$ws.Cells.Item(56, 1).Value = "Mallorca"
$ws.Cells.Item(56, 2).Value = 210
$ws.Cells.Item(56, 3).Value = 18
$ws.Cells.Item(56, 4).Value = 194
$ws.Cells.Item(56, 5).Value = 12
$ws.Cells.Item(57, 1).Value = "Huelva"
$ws.Cells.Item(57, 2).Value = 177
$ws.Cells.Item(57, 3).Value = 2
$ws.Cells.Item(57, 4).Value = 171
$ws.Cells.Item(57, 5).Value = 4
$ws.Cells.Item(58, 1).Value = "Melilla"
$ws.Cells.Item(58, 2).Value = 62
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 61
$ws.Cells.Item(58, 5).Value = 1
$ws.Cells.Item(59, 1).Value = "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena"
$ws.Cells.Item(59, 2).Value = 58
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 58
$ws.Cells.Item(59, 5).Value = 3
$ws.Cells.Item(60, 1).Value = "Ceuta"
$ws.Cells.Item(60, 2).Value = 29
$ws.Cells.Item(60, 3).Value = 0
$ws.Cells.Item(60, 4).Value = 28
$ws.Cells.Item(60, 5).Value = 1
$ws.Cells.Item(61, 1).Value = "Ibiza"
$ws.Cells.Item(61, 2).Value = 21
$ws.Cells.Item(61, 3).Value = 18
$ws.Cells.Item(61, 4).Value = 20
$ws.Cells.Item(61, 5).Value = 1
$ws.Cells.Item(62, 1).Value = "Menorca"
$ws.Cells.Item(62, 2).Value = 15
$ws.Cells.Item(62, 3).Value = 18
$ws.Cells.Item(62, 4).Value = 13
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 1).Value = "Arroyo de la Luz"
$ws.Cells.Item(63, 2).Value = 7
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 4).Value = 7
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 1).Value = "Formentera"
$ws.Cells.Item(64, 2).Value = 0
$ws.Cells.Item(64, 3).Value = 10
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 8
